$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 121.8
$ws.Range("I5").Value = 105.375
$ws.Range("J5").Value = 187.5
$ws.Range("K5").Value = 105.375
$ws.Range("L5").Value = 187.5
$ws.Range("M5").Value = 9.625
$ws.Range("N5").Value = -417.5
$ws.Range("H28").Value = 495.25
$ws.Range("I28").Value = 495.25
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 495.25
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -10.25
$ws.Range("H32").Value = 1998
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1998
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1998
$ws.Range("N32").Value = -2650
$ws.Range("H33").Value = 214.35484
$ws.Range("I33").Value = 183.46428
$ws.Range("J33").Value = 502.66666
$ws.Range("K33").Value = 183.46428
$ws.Range("L33").Value = 502.66666
$ws.Range("M33").Value = 45.53572
$ws.Range("H38").Value = 3218.6667
$ws.Range("I38").Value = 384
$ws.Range("J38").Value = 8888
$ws.Range("K38").Value = 1152
$ws.Range("L38").Value = 26664
$ws.Range("M38").Value = -780
$ws.Range("N38").Value = -27408
$ws.Range("H40").Value = 6999.769
$ws.Range("I40").Value = 11999.833
$ws.Range("J40").Value = 2714
$ws.Range("K40").Value = 11999.833
$ws.Range("L40").Value = 2714
$ws.Range("M40").Value = -11824.833
$ws.Range("N40").Value = -3064
$ws.Range("H42").Value = 672
$ws.Range("I42").Value = 311
$ws.Range("J42").Value = 1033
$ws.Range("K42").Value = 933
$ws.Range("L42").Value = 3099
$ws.Range("M42").Value = -703
$ws.Range("H76").Value = 8167.04
$ws.Range("I76").Value = 7722.8096
$ws.Range("J76").Value = 10499.25
$ws.Range("K76").Value = 7722.8096
$ws.Range("L76").Value = 10499.25
$ws.Range("M76").Value = -7407.8096
$ws.Range("H79").Value = 8167.04
$ws.Range("I79").Value = 7722.8096
$ws.Range("J79").Value = 10499.25
$ws.Range("K79").Value = 7722.8096
$ws.Range("L79").Value = 10499.25
$ws.Range("M79").Value = -6630.8096
$ws.Range("H100").Value = 1021
$ws.Range("I100").Value = 1021
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1021
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -480
$ws.Range("H116").Value = 3033.48
$ws.Range("I116").Value = 2864.9443
$ws.Range("J116").Value = 3466.8572
$ws.Range("K116").Value = 2864.9443
$ws.Range("L116").Value = 3466.8572
$ws.Range("M116").Value = 577.0556999999999
$ws.Range("N116").Value = -10350.8572
$ws.Range("H132").Value = 478104.38
$ws.Range("I132").Value = 1677.3334
$ws.Range("J132").Value = 3336666.8
$ws.Range("K132").Value = 5032.0002
$ws.Range("L132").Value = 10010000.4
$ws.Range("M132").Value = -2502.0002
$ws.Range("H137").Value = 3563.6562
$ws.Range("I137").Value = 1333.5
$ws.Range("J137").Value = 6431
$ws.Range("K137").Value = 4000.5
$ws.Range("L137").Value = 19293
$ws.Range("M137").Value = -1450.5
$ws.Range("H141").Value = 3893.8
$ws.Range("I141").Value = 2240.8
$ws.Range("J141").Value = 7199.8
$ws.Range("K141").Value = 6722.400000000001
$ws.Range("L141").Value = 21599.4
$ws.Range("M141").Value = -1542.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2350
$ws.Range("I2").Value = 2283.389
$ws.Range("J2").Value = 2949.5
$ws.Range("K2").Value = 2283.389
$ws.Range("L2").Value = 2949.5
$ws.Range("M2").Value = -2170.389
$ws.Range("H32").Value = 35688.543
$ws.Range("I32").Value = 10976.454
$ws.Range("J32").Value = 56598.77
$ws.Range("K32").Value = 10976.454
$ws.Range("L32").Value = 56598.77
$ws.Range("M32").Value = -10689.454
$ws.Range("H61").Value = 3459.2222
$ws.Range("I61").Value = 3039.6667
$ws.Range("J61").Value = 4298.3335
$ws.Range("K61").Value = 3039.6667
$ws.Range("L61").Value = 4298.3335
$ws.Range("M61").Value = -2827.6667
$ws.Range("N61").Value = -4722.3335
$ws.Range("H74").Value = 1029.5
$ws.Range("I74").Value = 1030.9048
$ws.Range("J74").Value = 1000
$ws.Range("K74").Value = 1030.9048
$ws.Range("L74").Value = 1000
$ws.Range("M74").Value = -156.9048
$ws.Range("H77").Value = 1029.5
$ws.Range("I77").Value = 1030.9048
$ws.Range("J77").Value = 1000
$ws.Range("K77").Value = 5154.524
$ws.Range("L77").Value = 5000
$ws.Range("M77").Value = -786.5240000000003
$ws.Range("H110").Value = 2048.2307
$ws.Range("I110").Value = 1662.1818
$ws.Range("J110").Value = 4171.5
$ws.Range("K110").Value = 1662.1818
$ws.Range("L110").Value = 4171.5
$ws.Range("M110").Value = 382.8181999999999
$ws.Range("H116").Value = 2350
$ws.Range("I116").Value = 2283.389
$ws.Range("J116").Value = 2949.5
$ws.Range("K116").Value = 2283.389
$ws.Range("L116").Value = 2949.5
$ws.Range("M116").Value = 10.61099999999988
$ws.Range("H132").Value = 3339.5881
$ws.Range("I132").Value = 1251.091
$ws.Range("J132").Value = 7168.5
$ws.Range("K132").Value = 3753.273
$ws.Range("L132").Value = 21505.5
$ws.Range("M132").Value = -1223.273
$ws.Range("N132").Value = -26565.5
$ws.Range("H134").Value = 67250
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 67250
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 67250
$ws.Range("N134").Value = -77390
$ws.Range("H136").Value = 3459.2222
$ws.Range("I136").Value = 3039.6667
$ws.Range("J136").Value = 4298.3335
$ws.Range("K136").Value = 9119.000100000001
$ws.Range("L136").Value = 12895.0005
$ws.Range("M136").Value = -6569.000100000001
$ws.Range("N136").Value = -17995.0005
$ws.Range("H139").Value = 61624.5
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 61624.5
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 61624.5
$ws.Range("N139").Value = -71904.5
$ws.Range("H141").Value = 49166.668
$ws.Range("I141").Value = 40000
$ws.Range("J141").Value = 67500
$ws.Range("K141").Value = 40000
$ws.Range("L141").Value = 67500
$ws.Range("M141").Value = -34820
$ws.Range("N141").Value = -77860

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2350
$ws.Range("I3").Value = 2283.389
$ws.Range("J3").Value = 2949.5
$ws.Range("K3").Value = 2283.389
$ws.Range("L3").Value = 2949.5
$ws.Range("M3").Value = -2169.389
$ws.Range("H80").Value = 594.2941
$ws.Range("I80").Value = 101.3
$ws.Range("J80").Value = 1298.5714
$ws.Range("K80").Value = 101.3
$ws.Range("L80").Value = 1298.5714
$ws.Range("M80").Value = 896.7
$ws.Range("H83").Value = 594.2941
$ws.Range("I83").Value = 101.3
$ws.Range("J83").Value = 1298.5714
$ws.Range("K83").Value = 506.5
$ws.Range("L83").Value = 6492.857
$ws.Range("M83").Value = 4485.5
$ws.Range("H99").Value = 2362.375
$ws.Range("I99").Value = 1580
$ws.Range("J99").Value = 3666.3333
$ws.Range("K99").Value = 1580
$ws.Range("L99").Value = 3666.3333
$ws.Range("M99").Value = -82
$ws.Range("H105").Value = 4337.0527
$ws.Range("I105").Value = 4491.5
$ws.Range("J105").Value = 4295.8667
$ws.Range("K105").Value = 4491.5
$ws.Range("L105").Value = 4295.8667
$ws.Range("M105").Value = -2744.5
$ws.Range("N105").Value = -7789.8667
$ws.Range("H134").Value = 2645.7144
$ws.Range("I134").Value = 2076.7334
$ws.Range("J134").Value = 4068.1667
$ws.Range("K134").Value = 6230.2002
$ws.Range("L134").Value = 12204.5001
$ws.Range("M134").Value = -3695.2002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 885.8889
$ws.Range("I22").Value = 799.75
$ws.Range("J22").Value = 954.8
$ws.Range("K22").Value = 799.75
$ws.Range("L22").Value = 954.8
$ws.Range("M22").Value = -449.75
$ws.Range("H31").Value = 4789.707
$ws.Range("I31").Value = 3233.1667
$ws.Range("J31").Value = 7336.773
$ws.Range("K31").Value = 3233.1667
$ws.Range("L31").Value = 7336.773
$ws.Range("M31").Value = -2938.1667
$ws.Range("H34").Value = 4789.707
$ws.Range("I34").Value = 3233.1667
$ws.Range("J34").Value = 7336.773
$ws.Range("K34").Value = 3233.1667
$ws.Range("L34").Value = 7336.773
$ws.Range("M34").Value = -3031.1667
$ws.Range("H52").Value = 87966
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 87966
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 87966
$ws.Range("N52").Value = -88554
$ws.Range("H69").Value = 44330.332
$ws.Range("I69").Value = 44330.332
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 44330.332
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -43581.332
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 44330.332
$ws.Range("I72").Value = 44330.332
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 132990.996
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -129246.996
$ws.Range("N72").ClearContents()
$ws.Range("H105").Value = 3253.6843
$ws.Range("I105").Value = 2755.6924
$ws.Range("J105").Value = 4332.6665
$ws.Range("K105").Value = 2755.6924
$ws.Range("L105").Value = 4332.6665
$ws.Range("M105").Value = -1008.6924
$ws.Range("H107").Value = 553.7646999999999
$ws.Range("I107").Value = 536.3077
$ws.Range("J107").Value = 610.5
$ws.Range("K107").Value = 536.3077
$ws.Range("L107").Value = 610.5
$ws.Range("M107").Value = 1383.6923
$ws.Range("H122").Value = 1742.3478
$ws.Range("I122").Value = 2033.1666
$ws.Range("J122").Value = 695.4
$ws.Range("K122").Value = 6099.4998
$ws.Range("L122").Value = 2086.2
$ws.Range("M122").Value = -3649.4998
$ws.Range("H132").Value = 4124.778
$ws.Range("I132").Value = 4395.6
$ws.Range("J132").Value = 3786.25
$ws.Range("K132").Value = 13186.8
$ws.Range("L132").Value = 11358.75
$ws.Range("M132").Value = -10656.8
$ws.Range("N132").Value = -16418.75
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H141").Value = 313959
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 313959
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 313959
$ws.Range("N141").Value = -324319

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 83.36842
$ws.Range("I2").Value = 59.857143
$ws.Range("J2").Value = 97.083336
$ws.Range("K2").Value = 359.142858
$ws.Range("L2").Value = 582.500016
$ws.Range("M2").Value = -246.142858
$ws.Range("N2").Value = -808.500016
$ws.Range("H5").Value = 698
$ws.Range("I5").Value = 698
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2094
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1982
$ws.Range("N5").ClearContents()
$ws.Range("H7").Value = 322.5
$ws.Range("I7").Value = 322.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 967.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -855.5
$ws.Range("H23").Value = 116.42857
$ws.Range("I23").Value = 132.5
$ws.Range("J23").Value = 20
$ws.Range("K23").Value = 397.5
$ws.Range("L23").Value = 60
$ws.Range("M23").Value = -162.5
$ws.Range("N23").Value = -530
$ws.Range("H34").Value = 343.81818
$ws.Range("I34").Value = 444.6
$ws.Range("J34").Value = 127.85714
$ws.Range("K34").Value = 1333.8
$ws.Range("L34").Value = 383.57142
$ws.Range("M34").Value = -1249.8
$ws.Range("N34").Value = -551.57142
$ws.Range("H46").Value = 158.11111
$ws.Range("I46").Value = 110.5
$ws.Range("J46").Value = 196.2
$ws.Range("K46").Value = 331.5
$ws.Range("L46").Value = 588.5999999999999
$ws.Range("M46").Value = -240.5
$ws.Range("N46").Value = -770.5999999999999
$ws.Range("H55").Value = 666889.7
$ws.Range("I55").Value = 270
$ws.Range("J55").Value = 1000199.5
$ws.Range("K55").Value = 810
$ws.Range("L55").Value = 3000598.5
$ws.Range("M55").Value = -633
$ws.Range("N55").Value = -3000952.5
$ws.Range("H104").Value = 342461.38
$ws.Range("I104").Value = 999
$ws.Range("J104").Value = 404545.47
$ws.Range("K104").Value = 2997
$ws.Range("L104").Value = 1213636.41
$ws.Range("M104").Value = -376
$ws.Range("N104").Value = -1218878.41
$ws.Range("H132").Value = 9038.526
$ws.Range("I132").Value = 15976
$ws.Range("J132").Value = 1330.2222
$ws.Range("K132").Value = 143784
$ws.Range("L132").Value = 11971.9998
$ws.Range("M132").Value = -141254
$ws.Range("N132").Value = -17031.9998
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("H135").Value = 698
$ws.Range("I135").Value = 698
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 6282
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -3747
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("H138").Value = 1794
$ws.Range("I138").Value = 1794
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 5382
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -242
$ws.Range("H139").Value = 23267.666
$ws.Range("I139").Value = 48999
$ws.Range("J139").Value = 18121.4
$ws.Range("K139").Value = 146997
$ws.Range("L139").Value = 54364.2
$ws.Range("M139").Value = -141857
$ws.Range("H140").Value = 2069.875
$ws.Range("I140").Value = 2015
$ws.Range("J140").Value = 2088.1667
$ws.Range("K140").Value = 6045
$ws.Range("L140").Value = 6264.500100000001
$ws.Range("M140").Value = -865
$ws.Range("N140").Value = -16624.5001
$ws.Range("H141").Value = 2000
$ws.Range("I141").Value = 2000
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 6000
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 54998.668
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 54998.668
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 54998.668
$ws.Range("N52").Value = -55516.668
$ws.Range("H80").Value = 8583.532999999999
$ws.Range("I80").Value = 4271.5884
$ws.Range("J80").Value = 14222.23
$ws.Range("K80").Value = 4271.5884
$ws.Range("L80").Value = 14222.23
$ws.Range("M80").Value = -3273.5884
$ws.Range("N80").Value = -16218.23
$ws.Range("H83").Value = 8583.532999999999
$ws.Range("I83").Value = 4271.5884
$ws.Range("J83").Value = 14222.23
$ws.Range("K83").Value = 21357.942
$ws.Range("L83").Value = 71111.14999999999
$ws.Range("M83").Value = -16365.942
$ws.Range("N83").Value = -81095.14999999999
$ws.Range("H107").Value = 678.64703
$ws.Range("I107").Value = 267.58334
$ws.Range("J107").Value = 1665.2
$ws.Range("K107").Value = 267.58334
$ws.Range("L107").Value = 1665.2
$ws.Range("M107").Value = 1652.41666
$ws.Range("N107").Value = -5505.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 344.64285
$ws.Range("I16").Value = 344.64285
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 344.64285
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -174.64285
$ws.Range("H22").Value = 855.4286
$ws.Range("H27").Value = 855.4286
$ws.Range("H50").Value = 35000
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 35000
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 35000
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -36274
$ws.Range("H53").Value = 5000
$ws.Range("I53").Value = 5000
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 5000
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -4482
$ws.Range("N53").ClearContents()
$ws.Range("H61").Value = 4001.125
$ws.Range("I61").Value = 4001.125
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4001.125
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3799.125
$ws.Range("H113").Value = 4001.125
$ws.Range("I113").Value = 4001.125
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4001.125
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1831.125
$ws.Range("H136").Value = 4099
$ws.Range("I136").Value = 3775.75
$ws.Range("J136").Value = 4530
$ws.Range("K136").Value = 11327.25
$ws.Range("L136").Value = 13590
$ws.Range("M136").Value = -8777.25
$ws.Range("N136").Value = -18690

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 3756.25
$ws.Range("I11").Value = 4006.6667
$ws.Range("J11").Value = 3005
$ws.Range("K11").Value = 4006.6667
$ws.Range("L11").Value = 3005
$ws.Range("M11").Value = -3864.6667
$ws.Range("N11").Value = -3289
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H29").Value = 33100
$ws.Range("I29").Value = 34650
$ws.Range("J29").Value = 30000
$ws.Range("K29").Value = 34650
$ws.Range("L29").Value = 30000
$ws.Range("M29").Value = -34360
$ws.Range("H34").Value = 43021.2
$ws.Range("I34").Value = 45035.668
$ws.Range("J34").Value = 39999.5
$ws.Range("K34").Value = 45035.668
$ws.Range("L34").Value = 39999.5
$ws.Range("M34").Value = -44832.668
$ws.Range("H107").Value = 1189.7059
$ws.Range("I107").Value = 684.6667
$ws.Range("J107").Value = 2401.8
$ws.Range("K107").Value = 2054.0001
$ws.Range("L107").Value = 7205.400000000001
$ws.Range("M107").Value = -134.0001000000002
$ws.Range("N107").Value = -11045.4
$ws.Range("H122").Value = 4452.1025
$ws.Range("I122").Value = 4053.0571
$ws.Range("J122").Value = 7943.75
$ws.Range("K122").Value = 12159.1713
$ws.Range("L122").Value = 23831.25
$ws.Range("M122").Value = -9709.1713
$ws.Range("H126").Value = 3200.6
$ws.Range("I126").Value = 3200.6
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9601.799999999999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7131.799999999999
$ws.Range("H132").Value = 1312.4667
$ws.Range("I132").Value = 1228.2858
$ws.Range("J132").Value = 2491
$ws.Range("K132").Value = 3684.8574
$ws.Range("L132").Value = 7473
$ws.Range("M132").Value = -1154.8574
$ws.Range("H136").Value = 4312.143
$ws.Range("I136").Value = 4312.143
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 12936.429
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -10386.429
$ws.Range("N136").ClearContents()
